# Update forecast summary values on the "Forecast Comparison" sheet.
# (Removed Auto Arima) - Amazon Mean/P70/P80/P90 Forecast columns (D:G)
# are recalculated for rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$newValues = @{
    2  = @(241, 287, 333, 404)
    3  = @(161, 194, 232, 291)
    4  = @(121, 147, 175, 220)
    5  = @(127, 154, 185, 233)
    6  = @(134, 163, 197, 252)
    7  = @(131, 160, 194, 250)
    8  = @(130, 159, 199, 266)
    9  = @(125, 152, 187, 243)
    10 = @(118, 143, 174, 223)
    11 = @(117, 142, 173, 222)
    12 = @(110, 134, 165, 214)
    13 = @(115, 141, 175, 231)
    14 = @(111, 135, 166, 215)
    15 = @(106, 129, 161, 213)
    16 = @(101, 123, 152, 200)
    17 = @(99, 121, 150, 196)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D: Amazon Mean Forecast
    $ws.Cells.Item($row, 5).Value = $vals[1]   # E: Amazon P70 Forecast
    $ws.Cells.Item($row, 6).Value = $vals[2]   # F: Amazon P80 Forecast
    $ws.Cells.Item($row, 7).Value = $vals[3]   # G: Amazon P90 Forecast
}
